# edit.ps1 - Applies the "minor changes to blog errors" commit to
# FinalReportSchedule.docx:
#   1. Splits the empty numbered (ilvl=0) bullet just above the final
#      "Implement Database Connectivity..." block into two paragraphs:
#      an empty indent-only paragraph, and a new "1.0" heading paragraph.
#   2. Moves the "_GoBack" bookmark from the end of the document into the
#      middle of the final "Implement Database Connectivity..." run
#      (splitting "between web pa" | "ge and database").
#   3. Appends two new bullet paragraphs after the final
#      "Incorporated error reporting..." paragraph:
#      "Major changes to system architecture:" (ilvl=1) and an empty
#      ilvl=2 bullet.

$d = $word.ActiveDocument

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$xmlFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Step 1: find the empty ilvl=0/numId=1 bullet that immediately precedes
# the final "Implement Database Connectivity..." bullet, and split it
# into two paragraphs (an empty indent-only paragraph and a new "1.0"
# paragraph).
# ---------------------------------------------------------------------
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Implement Database Connectivity")) {
        $targetIdx = $i
    }
}

$emptyIdx = $targetIdx - 1
$emptyPara = $d.Paragraphs.Item($emptyIdx)

$splitBody = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>1.0</w:t></w:r></w:p></w:body>'
$emptyPara.Range.InsertXML($xmlHeader + $splitBody + $xmlFooter)

# ---------------------------------------------------------------------
# Step 2: relocate the "_GoBack" bookmark into the middle of the final
# "Implement Database Connectivity..." run, splitting its text at
# "...between web pa|ge and database". Adding a bookmark with an
# existing name moves it (removing the old one automatically).
# ---------------------------------------------------------------------
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Implement Database Connectivity")) {
        $targetIdx = $i
    }
}

$connPara = $d.Paragraphs.Item($targetIdx)
$splitMarker = "Implement Database Connectivity between web pa"
$paraText = $connPara.Range.Text
$splitOffset = $paraText.IndexOf($splitMarker) + $splitMarker.Length
$splitPos = $connPara.Range.Start + $splitOffset
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# Step 3: append two new bullet paragraphs after the final
# "Incorporated error reporting..." paragraph.
# ---------------------------------------------------------------------
$lastIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Incorporated error reporting")) {
        $lastIdx = $i
    }
}

$lastPara = $d.Paragraphs.Item($lastIdx)
$endPos = $lastPara.Range.End
$endRange = $d.Range($endPos, $endPos)

$appendBody = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>Major changes to system architecture:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr></w:p></w:body>'
$endRange.InsertXML($xmlHeader + $appendBody + $xmlFooter)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
